$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.736.94"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "2.420.83"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.74"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.38"
$ws.Range("E6").Value = "  -3.15%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "2.405.63"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.06"
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.337"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.85"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("E15").Value = "  -4.22%  "
$ws.Range("D16").Value = "2.781.49"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").Value = "60.655.12"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "2.409.03"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.28"
$ws.Range("E19").Value = "  +12.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.61"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.09"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.81"
$ws.Range("E25").Value = "  -9.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.32"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "556.19"
$ws.Range("E27").Value = "  -6.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  -11.75%  "
$ws.Range("D30").Value = "0.0₃0916"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.31"
$ws.Range("E32").Value = "  -6.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  -4.73%  "
$ws.Range("E34").Value = "  -3.06%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.42"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.54"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.368"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.52"
$ws.Range("E39").Value = "  -5.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.24"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.11"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.65"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0290"
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.72"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.50"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.585"
$ws.Range("E48").Value = "  -3.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0498"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.04"
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0896"
$ws.Range("E51").Value = "  -0.40%  "
